$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: E2 "日期" -> "临床诊断病例", with new number-style (numFmt 0 / default font / bordered) ---
$ws.Range("E2").Value = "临床诊断病例"
$ws.Range("E2:E14").NumberFormat = "0"
$ws.Range("E2:E14").Borders.LineStyle = 1

# --- New "confirmed cases in hospital" column values (E3:E14) ---
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 4
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("E13").Value = 4
$ws.Range("E14").Value = 4

# --- Corrected count for 浠水 (row 8) ---
$ws.Range("B8").Value = 0

# --- Updated totals (row 13, literal values) ---
$ws.Range("B13").Value = 96

# --- Selection moves to B13 ---
$ws.Range("B13").Select()
